$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Include from MedComCorePracti").Name = "Include #0"
$wb.Worksheets.Item("Include from MedComCorePracti 2").Name = "Include #1"
$wb.Worksheets.Item("Include from NullFlavor").Name = "Include #2"

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.8.1"
